# Apply the AAPLNamed.xlsx update:
#  - Rename/insert header columns (ScoreFinal, totalSentiment, RSI, PEG)
#  - Rewrite the existing data row with new trading-run numbers
#  - Append a second data row for the new trading run
#  - Best-effort cosmetic tweaks (column widths, window size)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----------------------------------------------
$ws.Cells.Item(1, 1).Value2  = "Date"
$ws.Cells.Item(1, 2).Value2  = "ScoreFinal"
$ws.Cells.Item(1, 3).Value2  = "totalSentiment"
$ws.Cells.Item(1, 4).Value2  = "posWordPercentage"
$ws.Cells.Item(1, 5).Value2  = "negWordPercentage"
$ws.Cells.Item(1, 6).Value2  = "posPhrasePercentage"
$ws.Cells.Item(1, 7).Value2  = "negPhrasePercentage"
$ws.Cells.Item(1, 8).Value2  = "ElapsedMs"
$ws.Cells.Item(1, 9).Value2  = "wordCount"
$ws.Cells.Item(1, 10).Value2 = "sentenceCount"
$ws.Cells.Item(1, 11).Value2 = "posWordCount"
$ws.Cells.Item(1, 12).Value2 = "negWordCount"
$ws.Cells.Item(1, 13).Value2 = "positivePhraseCount"
$ws.Cells.Item(1, 14).Value2 = "negativePhraseCount"
$ws.Cells.Item(1, 15).Value2 = "Method"
$ws.Cells.Item(1, 16).Value2 = "RSI"
$ws.Cells.Item(1, 17).Value2 = "PEG"

# ---- Data row 2 (existing run, rewritten with new values) ------------
$ws.Cells.Item(2, 1).Value2  = 42627.874259259261
$ws.Cells.Item(2, 2).Value2  = 8
$ws.Cells.Item(2, 3).Value2  = 52
$ws.Cells.Item(2, 4).Value2  = 71
$ws.Cells.Item(2, 5).Value2  = 27
$ws.Cells.Item(2, 6).Value2  = 99
$ws.Cells.Item(2, 7).Value2  = 0
$ws.Cells.Item(2, 8).Value2  = 8756
$ws.Cells.Item(2, 9).Value2  = 5213
$ws.Cells.Item(2, 10).Value2 = 267
$ws.Cells.Item(2, 11).Value2 = 63
$ws.Cells.Item(2, 12).Value2 = 24
$ws.Cells.Item(2, 13).Value2 = 9
$ws.Cells.Item(2, 14).Value2 = 0
$ws.Cells.Item(2, 15).Value2 = "Named"
$ws.Cells.Item(2, 16).Value2 = 0
$ws.Cells.Item(2, 17).Value2 = 1

# ---- Data row 3 (new trading run) -------------------------------------
$ws.Cells.Item(3, 1).Value2  = 42627.877268518518
$ws.Cells.Item(3, 2).Value2  = 0
$ws.Cells.Item(3, 3).Value2  = -2
$ws.Cells.Item(3, 4).Value2  = 54
$ws.Cells.Item(3, 5).Value2  = 44
$ws.Cells.Item(3, 6).Value2  = 0
$ws.Cells.Item(3, 7).Value2  = 0
$ws.Cells.Item(3, 8).Value2  = 7368
$ws.Cells.Item(3, 9).Value2  = 3800
$ws.Cells.Item(3, 10).Value2 = 179
$ws.Cells.Item(3, 11).Value2 = 33
$ws.Cells.Item(3, 12).Value2 = 27
$ws.Cells.Item(3, 13).Value2 = 0
$ws.Cells.Item(3, 14).Value2 = 0
$ws.Cells.Item(3, 15).Value2 = "Named"
$ws.Cells.Item(3, 16).Value2 = 0
$ws.Cells.Item(3, 17).Value2 = 1

# ---- Column widths (best effort - engine quantizes to 1/6 char units) --
$ws.Columns.Item(1).ColumnWidth  = 14.0221354166667
$ws.Columns.Item(2).ColumnWidth  = 9.30729166666667
$ws.Columns.Item(3).ColumnWidth  = 13.7369791666667
$ws.Columns.Item(4).ColumnWidth  = 18.4518229166667
$ws.Columns.Item(5).ColumnWidth  = 18.5924479166667
$ws.Columns.Item(6).ColumnWidth  = 19.5924479166667
$ws.Columns.Item(7).ColumnWidth  = 19.7369791666667
$ws.Columns.Item(8).ColumnWidth  = 9.59244791666667
$ws.Columns.Item(9).ColumnWidth  = 10.0221354166667
$ws.Columns.Item(10).ColumnWidth = 13.7369791666667
$ws.Columns.Item(11).ColumnWidth = 13.5924479166667
$ws.Columns.Item(12).ColumnWidth = 13.7369791666667
$ws.Columns.Item(13).ColumnWidth = 18.8776041666667

# ---- Window size (best effort) ----------------------------------------
$win = $excel.ActiveWindow
$win.Width = 15105
$win.Height = 10215
